$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row
# (rows 2-189). The sheet was refreshed, bumping that date from
# 2023-10-04 (serial 45203) to 2023-10-05 (serial 45204) for every row.
for ($row = 2; $row -le 189; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
